$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Sheet1")

# ---------------------------------------------------------------------------
# Pump 1 block (rows 3-8 -> 3-9): a new "HP" parameter row is inserted just
# below the "Quantity" row. Only columns N:Q move down (the merged "Pump 1"
# header cell in column M, spanning M3:M7, and column M8 "***" stay put).
# ---------------------------------------------------------------------------

# Shift existing parameter rows down by one (process bottom-up so we never
# overwrite a source row before it has been copied).
$ws.Range("N8:Q8").Copy($ws.Range("N9:Q9"))
$ws.Range("N7:Q7").Copy($ws.Range("N8:Q8"))
$ws.Range("N6:Q6").Copy($ws.Range("N7:Q7"))
$ws.Range("N5:Q5").Copy($ws.Range("N6:Q6"))
$ws.Range("N4:Q4").Copy($ws.Range("N5:Q5"))

# Relabel the shifted rows per the new wording.
$ws.Range("N5").Value2 = "MaxGPM"
$ws.Range("Q6").Value2 = "Minimum pump turn down "

# Populate the brand-new row 4 with the "HP" parameter.
$ws.Range("N4").Value2 = "HP"
$ws.Range("P4").Value2 = "hp"
$ws.Range("Q4").Value2 = "design or break hoursepower"

# Row heights: the new Configuration row (9) needs the wrapped-text height
# that row 8 used to have; row 8's own (blank) row reverts to the default.
$ws.Rows.Item(9).RowHeight = 45

# ---------------------------------------------------------------------------
# Pump 2 block (rows 11-16 -> 12-17): the whole block (including column M,
# the merged "Pump 2" header) shifts down by one row so it stays aligned
# with the now-taller Pump 1 block above it.
# ---------------------------------------------------------------------------

$ws.Range("M16:Q16").Copy($ws.Range("M17:Q17"))
$ws.Range("M15:Q15").Copy($ws.Range("M16:Q16"))
$ws.Range("M14:Q14").Copy($ws.Range("M15:Q15"))
$ws.Range("M13:Q13").Copy($ws.Range("M14:Q14"))
$ws.Range("M12:Q12").Copy($ws.Range("M13:Q13"))

# Row 11's merged M cell must be unmerged before row 12 can receive its
# (now-separate) copy of the "Pump 2" header, and before row 11 is cleared.
$ws.Range("M11:M15").UnMerge()
$ws.Range("M11:Q11").Copy($ws.Range("M12:Q12"))

foreach ($col in @("M", "N", "O", "P", "Q")) {
    $ws.Range($col + "11").Clear()
}

$ws.Range("M12:M16").Merge()

# Row heights for the shifted rows.
$ws.Rows.Item(16).RowHeight = 45
$ws.Rows.Item(17).RowHeight = 45

# ---------------------------------------------------------------------------
# Cosmetic view state matching the saved workbook.
# ---------------------------------------------------------------------------
$ws.Application.ActiveWindow.ScrollColumn = 5
$ws.Range("M10").Select()
